# Update the PITONES price list: bump the sheet date one month forward and
# refresh all the unit prices (CON TOPE / SIN TOPE tables) with the new
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Header date (A1) moves from 2024-04-24 to 2024-05-24
$ws.Range("A1").Value = 45436

# CON TOPE prices (column D, rows 34-37)
$ws.Range("D34").Value = 2235.988
$ws.Range("D35").Value = 2981.315
$ws.Range("D36").Value = 4318.518
$ws.Range("D37").Value = 3672.706

# SIN TOPE prices (column D, rows 41-44)
$ws.Range("D41").Value = 2871.702
$ws.Range("D42").Value = 3794.088
$ws.Range("D43").Value = 4778.86
$ws.Range("D44").Value = 4176.859
